$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("screenerScenario3")
$ws2 = $wb.Worksheets.Item("IIQScenario3")
$ws3 = $wb.Worksheets.Item("RASSurveyScenario3")

# --- Data edits ---

# screenerScenario3
$ws1.Range("B5").Value = "ScreenerThreePLast"

$ws1.Range("A33").Value = "To determine eligibility for this study, we need to collect information about medical diagnoses."
$ws1.Rows.Item(33).RowHeight = 34

$ws1.Range("A37").Value = "We know that RASopathies are a group of syndromes that have a genetic basis.  In order to help us determine eligibility for the RASopathies Study, we also need to know about any genetic testing that has been completed."
$ws1.Rows.Item(37).RowHeight = 51

# IIQScenario3
$ws2.Range("B6").Value = "TestFirst TestLastThree"

# RASSurveyScenario3
$ws3.Range("B6").Value = "SurveyThreeFirst"
$ws3.Range("B8").Value = "SurveyThreeLast"

# --- View / selection changes ---
# (order matters: the last sheet selected/activated becomes the active tab,
#  so IIQScenario3 must be touched last to end up as the active sheet)

# screenerScenario3 is no longer the active tab; its stored selection moves to B13
$ws1.Range("B13").Select()

# RASSurveyScenario3 keeps its selection on B10 but loses its scrolled topLeftCell
$ws3.Range("B10").Select()

# IIQScenario3 becomes the active tab, selection on A3, no more scrolled topLeftCell
$ws2.Activate()
$ws2.Range("A3").Select()
